$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns (rows 2-51) remain text even for values that look numeric,
# matching the original inline-string / text cell type in the workbook.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '69.949.04'
$ws.Range('E2').Value = '  -0.54%  '

# Row 3
$ws.Range('D3').Value = '3.504.56'
$ws.Range('E3').Value = '  -1.69%  '

# Row 4
$ws.Range('E4').Value = '  -0.23%  '

# Row 5
$ws.Range('D5').Value = '607.69'
$ws.Range('E5').Value = '  -0.07%  '

# Row 6
$ws.Range('D6').Value = '197.46'
$ws.Range('E6').Value = '  +5.36%  '

# Row 7
$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +1.03%  '

# Row 8
$ws.Range('E8').Value = '  -0.12%  '

# Row 9
$ws.Range('D9').Value = '0.213'
$ws.Range('E9').Value = '  -0.19%  '

# Row 10
$ws.Range('D10').Value = '0.659'
$ws.Range('E10').Value = '  +1.84%  '

# Row 11
$ws.Range('D11').Value = '54.15'
$ws.Range('E11').Value = '  +0.38%  '

# Row 12
$ws.Range('D12').Value = '0.0000308'
$ws.Range('E12').Value = '  -0.33%  '

# Row 13
$ws.Range('D13').Value = '9.62'
$ws.Range('E13').Value = '  +2.36%  '

# Row 14
$ws.Range('D14').Value = '4.061.39'
$ws.Range('E14').Value = '  -1.76%  '

# Row 15
$ws.Range('D15').Value = '602.68'
$ws.Range('E15').Value = '  +5.17%  '

# Row 16
$ws.Range('D16').Value = '69.982.73'
$ws.Range('E16').Value = '  -0.60%  '

# Row 17
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '19.02'
$ws.Range('E17').Value = '  +0.26%  '

# Row 18
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '12.71'
$ws.Range('E18').Value = '  -0.28%  '

# Row 19
$ws.Range('D19').Value = '3.508.36'
$ws.Range('E19').Value = '  -1.89%  '

# Row 20
$ws.Range('E20').Value = '  -0.16%  '

# Row 21
$ws.Range('D21').Value = '0.994'
$ws.Range('E21').Value = '  +0.09%  '

# Row 22
$ws.Range('D22').Value = '17.87'
$ws.Range('E22').Value = '  +2.20%  '

# Row 23
$ws.Range('D23').Value = '104.64'
$ws.Range('E23').Value = '  +11.08%  '

# Row 24
$ws.Range('E24').Value = '  -2.43%  '

# Row 25
$ws.Range('D25').Value = '5.10'
$ws.Range('E25').Value = '  +2.89%  '

# Row 26
$ws.Range('E26').Value = '  +5.64%  '

# Row 27
$ws.Range('D27').Value = '10.99'
$ws.Range('E27').Value = '  +0.46%  '

# Row 28
$ws.Range('D28').Value = '9.85'
$ws.Range('E28').Value = '  +4.91%  '

# Row 29
$ws.Range('D29').Value = '34.02'
$ws.Range('E29').Value = '  +5.42%  '

# Row 30
$ws.Range('D30').Value = '4.60'
$ws.Range('E30').Value = '  +23.48%  '

# Row 31
$ws.Range('D31').Value = '7.22'
$ws.Range('E31').Value = '  +2.20%  '

# Row 32
$ws.Range('D32').Value = '12.69'
$ws.Range('E32').Value = '  +3.96%  '

# Row 33
$ws.Range('E33').Value = '  +1.11%  '

# Row 34
$ws.Range('D34').Value = '64.14'
$ws.Range('E34').Value = '  -0.39%  '

# Row 35
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = '524.09'
$ws.Range('E35').Value = '  +0.51%  '

# Row 36
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.687.12'
$ws.Range('E36').Value = '  -1.63%  '

# Row 37
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.15%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0800'
$ws.Range('E38').Value = '  +2.29%  '

# Row 39
$ws.Range('E39').Value = '  -4.95%  '

# Row 40
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').Value = '36.99'
$ws.Range('E40').Value = '  -2.36%  '

# Row 41
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.392'
$ws.Range('E41').Value = '  -3.01%  '

# Row 42
$ws.Range('E42').Value = '  +0.69%  '

# Row 43
$ws.Range('E43').Value = '  -1.16%  '

# Row 44
$ws.Range('E44').Value = '  +1.24%  '

# Row 45
$ws.Range('D45').Value = '2.87'
$ws.Range('E45').Value = '  -2.99%  '

# Row 46
$ws.Range('D46').Value = '0.141'
$ws.Range('E46').Value = '  +0.30%  '

# Row 47
$ws.Range('E47').Value = '  -4.71%  '

# Row 48
$ws.Range('D48').Value = '8.78'
$ws.Range('E48').Value = '  -4.42%  '

# Row 50
$ws.Range('D50').Value = '132.69'
$ws.Range('E50').Value = '  -1.65%  '

# Row 51
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = '0.000241'
$ws.Range('E51').Value = '  -1.61%  '

# Restore default (Normal) style so only the number format used for text-coercion
# is cleared from the affected cells, keeping their original appearance.
$textRange.Style = "Normal"
